# Apply odds/score updates per commit diff (Jogos_da_Semana_FlashScore_2025-03-16.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Range("AP28").Value = 3.15
$ws.Range("AQ28").Value = 1.36
$ws.Range("AR28").Value = 1.62
$ws.Range("AS28").Value = 2.31

# Row 29
$ws.Range("I29").Value = 5.5
$ws.Range("J29").Value = 2.4
$ws.Range("L29").Value = 6
$ws.Range("M29").Value = 1.1
$ws.Range("N29").Value = 7
$ws.Range("O29").Value = 1.44
$ws.Range("P29").Value = 2.63
$ws.Range("S29").Value = 4.5
$ws.Range("T29").Value = 1.18
$ws.Range("U29").Value = 1.53
$ws.Range("V29").Value = 2.38
$ws.Range("W29").Value = 2.2
$ws.Range("X29").Value = 1.62
$ws.Range("Z29").Value = 7
$ws.Range("AE29").Value = 7
$ws.Range("AG29").Value = 21
$ws.Range("AH29").Value = 81
$ws.Range("AJ29").Value = 11
$ws.Range("AK29").Value = 26
$ws.Range("AL29").Value = 19
$ws.Range("AN29").Value = 51
$ws.Range("AP29").Value = 3.7
$ws.Range("AQ29").Value = 1.28
$ws.Range("AR29").Value = 1.8
$ws.Range("AS29").Value = 2.05

# Row 30
$ws.Range("G30").Value = 2.1
$ws.Range("H30").Value = 3
$ws.Range("I30").Value = 4.1
$ws.Range("J30").Value = 3
$ws.Range("K30").Value = 1.83
$ws.Range("L30").Value = 5
$ws.Range("M30").Value = 1.14
$ws.Range("N30").Value = 5.5
$ws.Range("O30").Value = 1.57
$ws.Range("P30").Value = 2.25
$ws.Range("Q30").Value = 2.88
$ws.Range("R30").Value = 1.4
$ws.Range("S30").Value = 6
$ws.Range("T30").Value = 1.13
$ws.Range("U30").Value = 1.67
$ws.Range("V30").Value = 2.1
$ws.Range("Z30").Value = 8
$ws.Range("AA30").Value = 10
$ws.Range("AB30").Value = 19
$ws.Range("AC30").Value = 23
$ws.Range("AE30").Value = 5.5
$ws.Range("AF30").Value = 6
$ws.Range("AG30").Value = 21
$ws.Range("AJ30").Value = 8
$ws.Range("AK30").Value = 19
$ws.Range("AL30").Value = 15
$ws.Range("AM30").Value = 41
$ws.Range("AP30").Value = 4.8
$ws.Range("AQ30").Value = 1.19
$ws.Range("AR30").Value = 2.1
$ws.Range("AS30").Value = 1.78

# Row 37
$ws.Range("M37").Value = 1.2
$ws.Range("N37").Value = 4.33
$ws.Range("W37").Value = 3
$ws.Range("X37").Value = 1.36

# Row 64
$ws.Range("O64").Value = 1.18
$ws.Range("P64").Value = 4.5
$ws.Range("Q64").Value = 1.62
$ws.Range("R64").Value = 2.25
$ws.Range("S64").Value = 2.5
$ws.Range("T64").Value = 1.5
$ws.Range("AP64").Value = 2.02
$ws.Range("AQ64").Value = 1.77
$ws.Range("AR64").Value = 1.28
$ws.Range("AS64").Value = 3.55

# Row 67
$ws.Range("G67").Value = 1.53
$ws.Range("H67").Value = 3.8
$ws.Range("I67").Value = 7
$ws.Range("J67").Value = 2.1
$ws.Range("L67").Value = 7
$ws.Range("Y67").Value = 6
$ws.Range("Z67").Value = 6.5
$ws.Range("AB67").Value = 10
$ws.Range("AH67").Value = 67
$ws.Range("AK67").Value = 34
$ws.Range("AL67").Value = 21
$ws.Range("AO67").Value = 51

# Row 68
$ws.Range("G68").Value = 3.3
$ws.Range("AB68").Value = 41

# Row 86
$ws.Range("G86").Value = 3.7
$ws.Range("H86").Value = 3.1
$ws.Range("I86").Value = 2.15
$ws.Range("J86").Value = 4
$ws.Range("K86").Value = 2.05
$ws.Range("L86").Value = 2.88
$ws.Range("M86").Value = 1.08
$ws.Range("N86").Value = 8
$ws.Range("O86").Value = 1.36
$ws.Range("P86").Value = 3
$ws.Range("Q86").Value = 2.15
$ws.Range("R86").Value = 1.67
$ws.Range("W86").Value = 1.83
$ws.Range("X86").Value = 1.83
$ws.Range("Y86").Value = 10
$ws.Range("Z86").Value = 17
$ws.Range("AA86").Value = 13
$ws.Range("AB86").Value = 41
$ws.Range("AC86").Value = 29
$ws.Range("AD86").Value = 41
$ws.Range("AF86").Value = 6
$ws.Range("AG86").Value = 15
$ws.Range("AH86").Value = 51
$ws.Range("AI86").Value = 301
$ws.Range("AJ86").Value = 7
$ws.Range("AK86").Value = 10
$ws.Range("AL86").Value = 9.5
$ws.Range("AM86").Value = 19
$ws.Range("AN86").Value = 19
$ws.Range("AO86").Value = 29

# Row 129
$ws.Range("I129").Value = 7

# Row 132
$ws.Range("H132").Value = 5.25
$ws.Range("J132").Value = 12
$ws.Range("K132").Value = 2.5
$ws.Range("L132").Value = 1.67
$ws.Range("M132").Value = 1.06
$ws.Range("N132").Value = 10
$ws.Range("O132").Value = 1.25
$ws.Range("P132").Value = 3.75
$ws.Range("Q132").Value = 1.83
$ws.Range("R132").Value = 1.98
$ws.Range("S132").Value = 3
$ws.Range("T132").Value = 1.36
$ws.Range("U132").Value = 1.36
$ws.Range("V132").Value = 3
$ws.Range("W132").Value = 2.63
$ws.Range("X132").Value = 1.44
$ws.Range("Y132").Value = 21
$ws.Range("AA132").Value = 34
$ws.Range("AC132").Value = 101
$ws.Range("AD132").Value = 101
$ws.Range("AE132").Value = 10
$ws.Range("AH132").Value = 126
$ws.Range("AJ132").Value = 6
$ws.Range("AL132").Value = 10
$ws.Range("AN132").Value = 13

# Row 134
$ws.Range("G134").Value = 1.4
$ws.Range("H134").Value = 3.8
$ws.Range("M134").Value = 1.07
$ws.Range("N134").Value = 9
$ws.Range("O134").Value = 1.4
$ws.Range("P134").Value = 2.75
$ws.Range("Q134").Value = 2.2
$ws.Range("R134").Value = 1.65
$ws.Range("W134").Value = 2.5
$ws.Range("X134").Value = 1.5
$ws.Range("AB134").Value = 9
$ws.Range("AE134").Value = 7.5
$ws.Range("AG134").Value = 26

# Row 138
$ws.Range("G138").Value = 2.38
$ws.Range("I138").Value = 3.1
$ws.Range("J138").Value = 3.1
$ws.Range("L138").Value = 3.75
$ws.Range("Z138").Value = 11
$ws.Range("AA138").Value = 10
$ws.Range("AB138").Value = 23
$ws.Range("AJ138").Value = 8.5
$ws.Range("AN138").Value = 26

# Row 200
$ws.Range("H200").Value = 3.25
$ws.Range("J200").Value = 5.5
$ws.Range("K200").Value = 2.1
$ws.Range("M200").Value = 1.07
$ws.Range("N200").Value = 9
$ws.Range("O200").Value = 1.3
$ws.Range("P200").Value = 3.4
$ws.Range("Q200").Value = 2.05
$ws.Range("R200").Value = 1.8
$ws.Range("S200").Value = 3.5
$ws.Range("T200").Value = 1.29
$ws.Range("U200").Value = 1.44
$ws.Range("V200").Value = 2.63
$ws.Range("W200").Value = 1.91
$ws.Range("X200").Value = 1.8
$ws.Range("AE200").Value = 8.5
$ws.Range("AG200").Value = 17
$ws.Range("AI200").Value = 301
$ws.Range("AJ200").Value = 6.5
$ws.Range("AK200").Value = 7.5

# Row 203
$ws.Range("G203").Value = 2.7
$ws.Range("H203").Value = 3.9
$ws.Range("I203").Value = 2.3
$ws.Range("J203").Value = 3
$ws.Range("K203").Value = 2.6
$ws.Range("L203").Value = 2.75
$ws.Range("M203").Value = 1.01
$ws.Range("N203").Value = 23
$ws.Range("O203").Value = 1.08
$ws.Range("P203").Value = 8
$ws.Range("Q203").Value = 1.33
$ws.Range("R203").Value = 3.4
$ws.Range("S203").Value = 1.83
$ws.Range("T203").Value = 1.83
$ws.Range("Y203").Value = 17
$ws.Range("Z203").Value = 21
$ws.Range("AA203").Value = 11
$ws.Range("AB203").Value = 29
$ws.Range("AC203").Value = 19
$ws.Range("AE203").Value = 26
$ws.Range("AH203").Value = 23
$ws.Range("AJ203").Value = 17
$ws.Range("AN203").Value = 15
$ws.Range("AO203").Value = 17

# Row 206
$ws.Range("G206").Value = 1.85
$ws.Range("I206").Value = 4
$ws.Range("J206").Value = 2.37
$ws.Range("L206").Value = 4.3
$ws.Range("O206").Value = 1.31
$ws.Range("P206").Value = 2.9
$ws.Range("Q206").Value = 1.9
$ws.Range("S206").Value = 3.05
$ws.Range("W206").Value = 1.78
$ws.Range("X206").Value = 1.83
$ws.Range("Y206").Value = 6.8
$ws.Range("Z206").Value = 8.5
$ws.Range("AB206").Value = 15.5
$ws.Range("AD206").Value = 28
$ws.Range("AG206").Value = 15
$ws.Range("AH206").Value = 75
$ws.Range("AK206").Value = 22
$ws.Range("AL206").Value = 13.5
$ws.Range("AM206").Value = 65
$ws.Range("AN206").Value = 40
